$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Production (MW)" values for rows 2-97 (row 1 is header).
# Rows 43-97 keep their existing value of 0 (unchanged by this edit).
$newB = @{
    2 = 991;  3 = 0;    4 = 974;  5 = 987;  6 = 974;  7 = 975;  8 = 988;
    9 = 981;  10 = 968; 11 = 959; 12 = 946; 13 = 912; 14 = 884; 15 = 891;
    16 = 897; 17 = 887; 18 = 849; 19 = 816; 20 = 784; 21 = 778; 22 = 755;
    23 = 739; 24 = 692; 25 = 649; 26 = 584; 27 = 531; 28 = 539; 29 = 547;
    30 = 0;   31 = 0;   32 = 0;   33 = 0;   34 = 0;   35 = 0;   36 = 0;
    37 = 0;   38 = 0;   39 = 0;   40 = 0;   41 = 0;   42 = 0;   43 = 0;
    44 = 0;   45 = 0;   46 = 0;   47 = 0;   48 = 0;   49 = 0;   50 = 0;
    51 = 0;   52 = 0;   53 = 0;   54 = 0;   55 = 0;   56 = 0;   57 = 0;
    58 = 0;   59 = 0;   60 = 0;   61 = 0;   62 = 0;   63 = 0;   64 = 0;
    65 = 0;   66 = 0;   67 = 0;   68 = 0;   69 = 0;   70 = 0;   71 = 0;
    72 = 0;   73 = 0;   74 = 0;   75 = 0;   76 = 0;   77 = 0;   78 = 0;
    79 = 0;   80 = 0;   81 = 0;   82 = 0;   83 = 0;   84 = 0;   85 = 0;
    86 = 0;   87 = 0;   88 = 0;   89 = 0;   90 = 0;   91 = 0;   92 = 0;
    93 = 0;   94 = 0;   95 = 0;   96 = 0;   97 = 0
}

for ($r = 2; $r -le 97; $r++) {
    # Shift the timestamp forward by 11 days, keeping the same time-of-day.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value2 = $cellA.Value2 + 11

    # Apply the new production value.
    $ws.Cells.Item($r, 2).Value2 = $newB[$r]
}
